$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 91 and 92: the two match records were re-sequenced (their ids swapped
# places). Column A (the running "id" counter, 89 / 90) stays put; every
# other column (B..AB) from the two rows trades places.
# ---------------------------------------------------------------------------

# Row 91 <= data that used to live in row 92 (match id 6924569)
$ws.Cells.Item(91, 2).Value = 6924569
$ws.Cells.Item(91, 5).Value = "Venados FC"
$ws.Cells.Item(91, 6).Value = "Dorados"
$ws.Cells.Item(91, 7).Value = 4
$ws.Cells.Item(91, 8).Value = 1
$ws.Cells.Item(91, 9).Value = "H"
$ws.Cells.Item(91, 10).Value = 1.615
$ws.Cells.Item(91, 11).Value = 4
$ws.Cells.Item(91, 12).Value = 4.5
$ws.Cells.Item(91, 13).Value = 1.5
$ws.Cells.Item(91, 14).Value = 4.75
$ws.Cells.Item(91, 15).Value = 5.75
$ws.Cells.Item(91, 16).Value = -1.25
$ws.Cells.Item(91, 17).Value = 1.925
$ws.Cells.Item(91, 18).Value = 1.875
$ws.Cells.Item(91, 19).Value = 3
$ws.Cells.Item(91, 20).Value = 1.75
$ws.Cells.Item(91, 21).Value = 1.95
$ws.Cells.Item(91, 22).Value = 0.5
$ws.Cells.Item(91, 23).Value = -1
$ws.Cells.Item(91, 24).Value = -1
$ws.Cells.Item(91, 25).Value = 0.925
$ws.Cells.Item(91, 26).Value = -1
$ws.Cells.Item(91, 27).Value = 0.75
$ws.Cells.Item(91, 28).Value = -1

# Row 92 <= data that used to live in row 91 (match id 6924568)
$ws.Cells.Item(92, 2).Value = 6924568
$ws.Cells.Item(92, 5).Value = "Atletico Morelia"
$ws.Cells.Item(92, 6).Value = "Atlante"
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 1
$ws.Cells.Item(92, 9).Value = "A"
$ws.Cells.Item(92, 10).Value = 2.4
$ws.Cells.Item(92, 11).Value = 3
$ws.Cells.Item(92, 12).Value = 2.875
$ws.Cells.Item(92, 13).Value = 2.7
$ws.Cells.Item(92, 14).Value = 3.1
$ws.Cells.Item(92, 15).Value = 2.8
$ws.Cells.Item(92, 16).Value = 0
$ws.Cells.Item(92, 17).Value = 1.85
$ws.Cells.Item(92, 18).Value = 1.95
$ws.Cells.Item(92, 19).Value = 2.25
$ws.Cells.Item(92, 20).Value = 1.975
$ws.Cells.Item(92, 21).Value = 1.725
$ws.Cells.Item(92, 22).Value = -1
$ws.Cells.Item(92, 23).Value = -1
$ws.Cells.Item(92, 24).Value = 1.8
$ws.Cells.Item(92, 25).Value = -1
$ws.Cells.Item(92, 26).Value = 0.95
$ws.Cells.Item(92, 27).Value = -1
$ws.Cells.Item(92, 28).Value = 0.7250000000000001

# ---------------------------------------------------------------------------
# Rows 186 and 187: same kind of re-sequencing between two match records.
# ---------------------------------------------------------------------------

# Row 186 <= data that used to live in row 187 (match id 7648958)
$ws.Cells.Item(186, 2).Value = 7648958
$ws.Cells.Item(186, 5).Value = "Monterrey U23"
$ws.Cells.Item(186, 6).Value = "Mazatlan FC U23"
$ws.Cells.Item(186, 7).Value = 4
$ws.Cells.Item(186, 8).Value = 3
$ws.Cells.Item(186, 9).Value = "H"
$ws.Cells.Item(186, 10).Value = 2.375
$ws.Cells.Item(186, 11).Value = 3.1
$ws.Cells.Item(186, 12).Value = 2.75
$ws.Cells.Item(186, 13).Value = 2.375
$ws.Cells.Item(186, 14).Value = 3.4
$ws.Cells.Item(186, 15).Value = 3
$ws.Cells.Item(186, 16).Value = -0.25
$ws.Cells.Item(186, 17).Value = 2
$ws.Cells.Item(186, 18).Value = 1.8
$ws.Cells.Item(186, 19).Value = 2.75
$ws.Cells.Item(186, 20).Value = 1.95
$ws.Cells.Item(186, 21).Value = 1.85
$ws.Cells.Item(186, 22).Value = 1.375
$ws.Cells.Item(186, 23).Value = -1
$ws.Cells.Item(186, 24).Value = -1
$ws.Cells.Item(186, 25).Value = 1
$ws.Cells.Item(186, 26).Value = -1
$ws.Cells.Item(186, 27).Value = 0.95
$ws.Cells.Item(186, 28).Value = -1

# Row 187 <= data that used to live in row 186 (match id 7648957)
$ws.Cells.Item(187, 2).Value = 7648957
$ws.Cells.Item(187, 5).Value = "Unam Pumas U23"
$ws.Cells.Item(187, 6).Value = "Tijuana U23"
$ws.Cells.Item(187, 7).Value = 2
$ws.Cells.Item(187, 8).Value = 0
$ws.Cells.Item(187, 9).Value = "H"
$ws.Cells.Item(187, 10).Value = 1.666
$ws.Cells.Item(187, 11).Value = 3.5
$ws.Cells.Item(187, 12).Value = 4.2
$ws.Cells.Item(187, 13).Value = 1.533
$ws.Cells.Item(187, 14).Value = 4.333
$ws.Cells.Item(187, 15).Value = 6
$ws.Cells.Item(187, 16).Value = -1.25
$ws.Cells.Item(187, 17).Value = 2.025
$ws.Cells.Item(187, 18).Value = 1.775
$ws.Cells.Item(187, 19).Value = 2.75
$ws.Cells.Item(187, 20).Value = 1.775
$ws.Cells.Item(187, 21).Value = 2.025
$ws.Cells.Item(187, 22).Value = 0.5329999999999999
$ws.Cells.Item(187, 23).Value = -1
$ws.Cells.Item(187, 24).Value = -1
$ws.Cells.Item(187, 25).Value = 1.025
$ws.Cells.Item(187, 26).Value = -1
$ws.Cells.Item(187, 27).Value = -1
$ws.Cells.Item(187, 28).Value = 1.025

# ---------------------------------------------------------------------------
# New match appended as row 249 (id 247). Copy the formatting of the last
# existing data row first (bold/bordered/centered id cell in column A, and
# the date-time number format in column D), then fill in the values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(248, 1).Copy($ws.Cells.Item(249, 1))
$ws.Cells.Item(248, 4).Copy($ws.Cells.Item(249, 4))

$ws.Cells.Item(249, 1).Value = 247
$ws.Cells.Item(249, 2).Value = 8219587
$ws.Cells.Item(249, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(249, 4).Value = 45427.91666666666
$ws.Cells.Item(249, 5).Value = "Cancun FC"
$ws.Cells.Item(249, 6).Value = "Atlante"
$ws.Cells.Item(249, 7).Value = 1
$ws.Cells.Item(249, 8).Value = 1
$ws.Cells.Item(249, 9).Value = "D"
$ws.Cells.Item(249, 10).Value = 2.5
$ws.Cells.Item(249, 11).Value = 3
$ws.Cells.Item(249, 12).Value = 2.6
$ws.Cells.Item(249, 13).Value = 3.1
$ws.Cells.Item(249, 14).Value = 3
$ws.Cells.Item(249, 15).Value = 2.3
$ws.Cells.Item(249, 16).Value = 0.25
$ws.Cells.Item(249, 17).Value = 1.8
$ws.Cells.Item(249, 18).Value = 2
$ws.Cells.Item(249, 19).Value = 1.75
$ws.Cells.Item(249, 20).Value = 1.8
$ws.Cells.Item(249, 21).Value = 2
$ws.Cells.Item(249, 22).Value = -1
$ws.Cells.Item(249, 23).Value = 2
$ws.Cells.Item(249, 24).Value = -1
$ws.Cells.Item(249, 25).Value = 0.4
$ws.Cells.Item(249, 26).Value = -0.5
$ws.Cells.Item(249, 27).Value = 0.4
$ws.Cells.Item(249, 28).Value = -0.5
